$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.449754000000001
$ws.Range("H2").Value = 19.349262
$ws.Range("I2").Value = 0.03479900749229446
$ws.Range("J2").Value = 0.03479900749229446
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 4.798358985840001
$ws.Range("R2").Value = 43.18523087256001
$ws.Range("S2").Value = 0.0002178070998358664
$ws.Range("T2").Value = 0.0002178070998358664
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.449754000000001
$ws.Range("H3").Value = 19.349262
$ws.Range("I3").Value = 0.03479900749229446
$ws.Range("J3").Value = 0.03479900749229446
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 568.5444305528081
$ws.Range("R3").Value = 5116.899874975272
$ws.Range("S3").Value = 0.02580736745874445
$ws.Range("T3").Value = 0.02580736745874445
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.449754000000001
$ws.Range("H4").Value = 19.349262
$ws.Range("I4").Value = 0.03479900749229446
$ws.Range("J4").Value = 0.03479900749229446
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 192.000142624564
$ws.Range("R4").Value = 1728.001283621076
$ws.Range("S4").Value = 0.008715270023884669
$ws.Range("T4").Value = 0.008715270023884669
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.449754000000001
$ws.Range("H5").Value = 19.349262
$ws.Range("I5").Value = 0.03479900749229446
$ws.Range("J5").Value = 0.03479900749229446
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 1.290159342046
$ws.Range("R5").Value = 11.611434078414
$ws.Range("S5").Value = 0.00005856290982947288
$ws.Range("T5").Value = 0.00005856290982947289
$ws.Range("I6").Value = 0.663783921437469
$ws.Range("J6").Value = 0.6637839214374691
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 91.52771224268
$ws.Range("R6").Value = 823.7494101841199
$ws.Range("S6").Value = 0.004154625699540061
$ws.Range("T6").Value = 0.004154625699540062
$ws.Range("I7").Value = 0.663783921437469
$ws.Range("J7").Value = 0.6637839214374691
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("R7").Value = 97603.81428596983
$ws.Range("S7").Value = 0.4922702343604578
$ws.Range("T7").Value = 0.4922702343604579
$ws.Range("I8").Value = 0.663783921437469
$ws.Range("J8").Value = 0.6637839214374691
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 3662.363290565311
$ws.Range("R8").Value = 32961.2696150878
$ws.Range("S8").Value = 0.166241985899212
$ws.Range("T8").Value = 0.1662419858992121
$ws.Range("I9").Value = 0.663783921437469
$ws.Range("J9").Value = 0.6637839214374691
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 24.60952449670034
$ws.Range("R9").Value = 221.485720470303
$ws.Range("S9").Value = 0.001117075478259088
$ws.Range("T9").Value = 0.001117075478259088
$ws.Range("G10").Value = 55.79038633333334
$ws.Range("H10").Value = 167.371159
$ws.Range("I10").Value = 0.3010114916028843
$ws.Range("J10").Value = 0.3010114916028843
$ws.Range("M10").Value = 0.74396
$ws.Range("N10").Value = 2.23188
$ws.Range("O10").Value = 0.006259003216804254
$ws.Range("P10").Value = 0.006259003216804255
$ws.Range("Q10").Value = 41.50581581654667
$ws.Range("R10").Value = 373.55234234892
$ws.Range("S10").Value = 0.001884031894237499
$ws.Range("T10").Value = 0.0018840318942375
$ws.Range("G11").Value = 55.79038633333334
$ws.Range("H11").Value = 167.371159
$ws.Range("I11").Value = 0.3010114916028843
$ws.Range("J11").Value = 0.3010114916028843
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("O11").Value = 0.7416121699579786
$ws.Range("P11").Value = 0.7416121699579786
$ws.Range("Q11").Value = 4917.910578947067
$ws.Range("R11").Value = 44261.1952105236
$ws.Range("S11").Value = 0.2232337854699029
$ws.Range("T11").Value = 0.2232337854699029
$ws.Range("G12").Value = 55.79038633333334
$ws.Range("H12").Value = 167.371159
$ws.Range("I12").Value = 0.3010114916028843
$ws.Range("J12").Value = 0.3010114916028843
$ws.Range("M12").Value = 29.76859933333333
$ws.Range("N12").Value = 89.305798
$ws.Range("O12").Value = 0.2504459365921425
$ws.Range("P12").Value = 0.2504459365921425
$ws.Range("Q12").Value = 1660.801657408876
$ws.Range("R12").Value = 14947.21491667988
$ws.Range("S12").Value = 0.07538710493948218
$ws.Range("T12").Value = 0.07538710493948218
$ws.Range("G13").Value = 55.79038633333334
$ws.Range("H13").Value = 167.371159
$ws.Range("I13").Value = 0.3010114916028843
$ws.Range("J13").Value = 0.3010114916028843
$ws.Range("M13").Value = 0.2000323333333334
$ws.Range("N13").Value = 0.6000970000000001
$ws.Range("O13").Value = 0.00168289023307462
$ws.Range("P13").Value = 0.00168289023307462
$ws.Range("Q13").Value = 11.15988115582478
$ws.Range("R13").Value = 100.438930402423
$ws.Range("S13").Value = 0.000506569299261717
$ws.Range("T13").Value = 0.0005065692992617171
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.07517133333333333
$ws.Range("H14").Value = 0.225514
$ws.Range("I14").Value = 0.0004055794673521549
$ws.Range("J14").Value = 0.000405579467352155
$ws.Range("M14").Value = 0.74396
$ws.Range("N14").Value = 2.23188
$ws.Range("O14").Value = 0.006259003216804254
$ws.Range("P14").Value = 0.006259003216804255
$ws.Range("Q14").Value = 0.05592446514666666
$ws.Range("R14").Value = 0.50332018632
$ws.Range("S14").Value = 0.000002538523190826893
$ws.Range("T14").Value = 0.000002538523190826894
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.07517133333333333
$ws.Range("H15").Value = 0.225514
$ws.Range("I15").Value = 0.0004055794673521549
$ws.Range("J15").Value = 0.000405579467352155
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("O15").Value = 0.7416121699579786
$ws.Range("P15").Value = 0.7416121699579786
$ws.Range("Q15").Value = 6.626336896553776
$ws.Range("R15").Value = 59.63703206898399
$ws.Range("S15").Value = 0.0003007826688734328
$ws.Range("T15").Value = 0.0003007826688734328
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.07517133333333333
$ws.Range("H16").Value = 0.225514
$ws.Range("I16").Value = 0.0004055794673521549
$ws.Range("J16").Value = 0.000405579467352155
$ws.Range("M16").Value = 29.76859933333333
$ws.Range("N16").Value = 89.305798
$ws.Range("O16").Value = 0.2504459365921425
$ws.Range("P16").Value = 0.2504459365921425
$ws.Range("Q16").Value = 2.237745303352444
$ws.Range("R16").Value = 20.139707730172
$ws.Range("S16").Value = 0.0001015757295635527
$ws.Range("T16").Value = 0.0001015757295635527
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.07517133333333333
$ws.Range("H17").Value = 0.225514
$ws.Range("I17").Value = 0.0004055794673521549
$ws.Range("J17").Value = 0.000405579467352155
$ws.Range("M17").Value = 0.2000323333333334
$ws.Range("N17").Value = 0.6000970000000001
$ws.Range("O17").Value = 0.00168289023307462
$ws.Range("P17").Value = 0.00168289023307462
$ws.Range("Q17").Value = 0.01503669720644444
$ws.Range("R17").Value = 0.135330274858
$ws.Range("S17").Value = 0.0000006825457243425483
$ws.Range("T17").Value = 0.0000006825457243425485
